$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Column C ("Förändrad") for rows 2..507 moves from 2023-10-06 (45205) to 2023-10-07 (45206)
for ($r = 2; $r -le 507; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}

# 2) Row 507 gains an explicit row height (matches the rest of the sheet)
$ws.Rows.Item(507).RowHeight = 15

# 3) Append new row 508 with the new cleared-notification record
$newRow = 508

# Copy number formats from the row above so the new date cells (B/C) and the
# wrapped "Artnamn" cell (R) keep the same formatting as every other row.
$ws.Cells.Item($newRow, 2).NumberFormat = $ws.Cells.Item($newRow - 1, 2).NumberFormat
$ws.Cells.Item($newRow, 3).NumberFormat = $ws.Cells.Item($newRow - 1, 3).NumberFormat
$ws.Cells.Item($newRow, 18).WrapText = $ws.Cells.Item($newRow - 1, 18).WrapText

$ws.Cells.Item($newRow, 1).Value = "A 48029-2023"
$ws.Cells.Item($newRow, 2).Value = 45204
$ws.Cells.Item($newRow, 3).Value = 45206
$ws.Cells.Item($newRow, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item($newRow, 5).Value = "LEKSAND"
$ws.Cells.Item($newRow, 7).Value = 2.1
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0
